$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CreateNewCustomer")
$ws.Range("H2").Value = "n"
